$wb = $excel.ActiveWorkbook

# --- Football sheet: C3:C14 -> centered horizontal alignment (new style) ---
$wsFootball = $wb.Worksheets.Item("Football")
$rngFootball = $wsFootball.Range("C3:C14")
$rngFootball.HorizontalAlignment = -4108  # xlCenter

# --- Volleyball sheet: C3:C6 -> centered horizontal alignment (new style) ---
$wsVolleyball = $wb.Worksheets.Item("Volleyball")
$rngVolleyball = $wsVolleyball.Range("C3:C6")
$rngVolleyball.HorizontalAlignment = -4108  # xlCenter

# --- Basketball sheet: C3:C5 -> centered horizontal alignment (new style) + new values ---
$wsBasketball = $wb.Worksheets.Item("Basketball")
$rngBasketball = $wsBasketball.Range("C3:C5")
$rngBasketball.HorizontalAlignment = -4108  # xlCenter

$wsBasketball.Range("C3").Value = 0.2
$wsBasketball.Range("C4").Value = 0.2
$wsBasketball.Range("C5").Value = 0.2
